$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Range("H31").Value = 603
$ws.Range("I31").Value = 603
$ws.Range("K31").Value = 1809
$ws.Range("M31").Value = -1579

# Row 113
$ws.Range("H113").Value = 7694410.5
$ws.Range("I113").Value = 11112888
$ws.Range("K113").Value = 11112888
$ws.Range("M113").Value = -11109634

# Row 137
$ws.Range("H137").Value = 17872294
$ws.Range("I137").Value = 1117.04
$ws.Range("J137").Value = 62550240
$ws.Range("K137").Value = 3351.12
$ws.Range("L137").Value = 187650720
$ws.Range("M137").Value = -801.1199999999999
$ws.Range("N137").Value = -187655820

# Row 141
$ws.Range("H141").Value = 1488.7142
$ws.Range("I141").Value = 1135.8928
$ws.Range("J141").Value = 2900
$ws.Range("K141").Value = 3407.6784
$ws.Range("L141").Value = 8700
$ws.Range("M141").Value = 1772.3216
$ws.Range("N141").Value = -19060


$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = 0

# Row 6
$ws.Range("H6").Value = 300
$ws.Range("I6").Value = 300
$ws.Range("K6").Value = 300
$ws.Range("M6").Value = -127

# Row 22
$ws.Range("H22").Value = 3252.8
$ws.Range("I22").Value = 3252.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3252.8
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2953.8

# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

# Row 32
$ws.Range("H32").Value = 4259.15
$ws.Range("I32").Value = 3930.8523
$ws.Range("J32").Value = 6666.6665
$ws.Range("K32").Value = 3930.8523
$ws.Range("L32").Value = 6666.6665
$ws.Range("M32").Value = -3643.8523
$ws.Range("N32").Value = -7240.6665

# Row 132
$ws.Range("H132").Value = 11113832
$ws.Range("I132").Value = 13161246
$ws.Range("J132").Value = 4630356
$ws.Range("K132").Value = 39483738
$ws.Range("L132").Value = 13891068
$ws.Range("M132").Value = -39481208
$ws.Range("N132").Value = -13896128


$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = 0

# Row 107
$ws.Range("H107").Value = 417342.12
$ws.Range("I107").Value = 910027.6
$ws.Range("J107").Value = 454.3846
$ws.Range("K107").Value = 910027.6
$ws.Range("L107").Value = 454.3846
$ws.Range("M107").Value = -908107.6
$ws.Range("N107").Value = -4294.3846

# Row 132
$ws.Range("H132").Value = 39540
$ws.Range("J132").Value = 39540
$ws.Range("L132").Value = 39540
$ws.Range("N132").Value = -49660


$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 90.72727
$ws.Range("I7").Value = 66.333336
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 66.333336
$ws.Range("L7").Value = 120
$ws.Range("M7").Value = 46.666664
$ws.Range("N7").Value = -346

# Row 16
$ws.Range("H16").Value = 1779.95
$ws.Range("I16").Value = 1723
$ws.Range("J16").Value = 1885.7142
$ws.Range("K16").Value = 1723
$ws.Range("L16").Value = 1885.7142
$ws.Range("M16").Value = -1436
$ws.Range("N16").Value = -2459.7142

# Row 22
$ws.Range("H22").Value = 433.46667
$ws.Range("I22").Value = 191.09091
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 191.09091
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = 158.90909
$ws.Range("N22").Value = -1800

# Row 31
$ws.Range("H31").Value = 1160039.5
$ws.Range("I31").Value = 1027.2683
$ws.Range("J31").Value = 4815385.5
$ws.Range("K31").Value = 1027.2683
$ws.Range("L31").Value = 4815385.5
$ws.Range("M31").Value = -732.2683
$ws.Range("N31").Value = -4815975.5

# Row 34
$ws.Range("H34").Value = 1160039.5
$ws.Range("I34").Value = 1027.2683
$ws.Range("J34").Value = 4815385.5
$ws.Range("K34").Value = 1027.2683
$ws.Range("L34").Value = 4815385.5
$ws.Range("M34").Value = -825.2683
$ws.Range("N34").Value = -4815789.5

# Row 113
$ws.Range("H113").Value = 1779.95
$ws.Range("I113").Value = 1723
$ws.Range("J113").Value = 1885.7142
$ws.Range("K113").Value = 1723
$ws.Range("L113").Value = 1885.7142
$ws.Range("M113").Value = 447
$ws.Range("N113").Value = -6225.7142

# Row 134
$ws.Range("H134").Value = 889869.3
$ws.Range("I134").Value = 978.4878
$ws.Range("J134").Value = 10001000
$ws.Range("K134").Value = 2935.4634
$ws.Range("L134").Value = 30003000
$ws.Range("M134").Value = -400.4634000000001
$ws.Range("N134").Value = -30008070


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 50739.7
$ws.Range("I16").Value = 59523.176
$ws.Range("K16").Value = 59523.176
$ws.Range("M16").Value = -59353.176

# Row 35
$ws.Range("H35").Value = 298.57144
$ws.Range("I35").Value = 298.57144
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 298.57144
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = 37.42856

# Row 68
$ws.Range("H68").Value = 2363.1428
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 2508.4
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 2508.4
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4006.4

# Row 71
$ws.Range("H71").Value = 2363.1428
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 2508.4
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 12542
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -20030

# Row 100
$ws.Range("H100").Value = 3306.7742
$ws.Range("I100").Value = 1200.9286
$ws.Range("J100").Value = 5041
$ws.Range("K100").Value = 1200.9286
$ws.Range("L100").Value = 5041
$ws.Range("M100").Value = -659.9286
$ws.Range("N100").Value = -6123

# Row 122
$ws.Range("H122").Value = 10538908
$ws.Range("I122").Value = 1331823
$ws.Range("J122").Value = 40001580
$ws.Range("K122").Value = 3995469
$ws.Range("L122").Value = 120004740
$ws.Range("M122").Value = -3993019
$ws.Range("N122").Value = -120009640

# Row 132
$ws.Range("H132").Value = 3175752.5
$ws.Range("I132").Value = 3402277.8
$ws.Range("J132").Value = 4400
$ws.Range("K132").Value = 10206833.4
$ws.Range("L132").Value = 13200
$ws.Range("M132").Value = -10204303.4
$ws.Range("N132").Value = -18260


$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 26855.666
$ws.Range("J15").Value = 26855.666
$ws.Range("L15").Value = 26855.666
$ws.Range("N15").Value = -27431.666

